$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Refresh the scrape timestamp (column O) for every data row (2..73)
# ---------------------------------------------------------------------
$ws.Range("O2:O73").Value = "2022-08-23 20:57:49"

# ---------------------------------------------------------------------
# 2) The crawler re-ordered the scraped products - rows 10-13 get
#    cyclically shifted: old row 11 -> row 10, old row 12 -> row 11,
#    old row 13 -> row 12, old row 10 -> row 13 (columns A-N; column O
#    has already been refreshed with the new timestamp above).
#    Columns, in order, are:
#    A id, B title, C href, D quantity, E ratingAmount, F ratingValue,
#    G brand, H price, I priceContext, J priceContextHiddenText,
#    K priceContextPrice, L priceContextAmount, M udoCat, N productAriaLabel
# ---------------------------------------------------------------------

$newRow10 = "6075745006", "Avela Strumpfhose Madame Hasel  9.5", "/de/haushalt-tier/bekleidung/socken-unterwaesche/struempfe/avela-strumpfhose-madame-hasel-95/p/6075745006", $null, $null, 0, "Avela", "5.95", $null, $null, $null, $null, "['haushalt-tier', 'bekleidung', 'socken-unterwaesche', 'struempfe']", "Avela Strumpfhose Madame Hasel  9.5 5.95 Schweizer Franken"
$newRow11 = "3404677005", "Naturaline Herren Slip schwarz S", "/de/haushalt-tier/bekleidung/socken-unterwaesche/unterwaesche/naturaline-herren-slip-schwarz-s/p/3404677005", $null, $null, 0, "Coop", "14.95", $null, $null, $null, $null, "['haushalt-tier', 'bekleidung', 'socken-unterwaesche', 'unterwaesche']", "Naturaline Herren Slip schwarz S 14.95 Schweizer Franken"
$newRow12 = "6031467010", "Naturaline Herren T-Shirt Kurzarm weiss XXL", "/de/haushalt-tier/bekleidung/shirts-pullover/herren-shirt/naturaline-herren-t-shirt-kurzarm-weiss-xxl/p/6031467010", $null, $null, 0, "Coop", "24.95", $null, $null, $null, $null, "['haushalt-tier', 'bekleidung', 'shirts-pullover', 'herren-shirt']", "Naturaline Herren T-Shirt Kurzarm weiss XXL 24.95 Schweizer Franken"
$newRow13 = "6866452", "Selenacare Menstruationsunterwäsche M", "/de/haushalt-tier/bekleidung/socken-unterwaesche/unterwaesche/selenacare-menstruationsunterwaesche-m/p/6866452", "1ST", $null, 0, "Selenacare", "24.50", "24.50/1ST", "Preis pro 1 Stück", "24.50", "1ST", "['haushalt-tier', 'bekleidung', 'socken-unterwaesche', 'unterwaesche']", "Selenacare Menstruationsunterwäsche M 24.50 Schweizer Franken"

# Row 10 --------------------------------------------------------------
$ws.Cells.Item(10, 1).NumberFormat = "@"
$ws.Cells.Item(10, 1).Value = $newRow10[0]
$ws.Cells.Item(10, 2).Value = $newRow10[1]
$ws.Cells.Item(10, 3).Value = $newRow10[2]
$ws.Cells.Item(10, 4).Value = $newRow10[3]
$ws.Cells.Item(10, 5).Value = $newRow10[4]
$ws.Cells.Item(10, 6).Value = $newRow10[5]
$ws.Cells.Item(10, 7).Value = $newRow10[6]
$ws.Cells.Item(10, 8).NumberFormat = "@"
$ws.Cells.Item(10, 8).Value = $newRow10[7]
$ws.Cells.Item(10, 9).Value = $newRow10[8]
$ws.Cells.Item(10, 10).Value = $newRow10[9]
$ws.Cells.Item(10, 11).Value = $newRow10[10]
$ws.Cells.Item(10, 12).Value = $newRow10[11]
$ws.Cells.Item(10, 13).Value = $newRow10[12]
$ws.Cells.Item(10, 14).Value = $newRow10[13]

# Row 11 --------------------------------------------------------------
$ws.Cells.Item(11, 1).NumberFormat = "@"
$ws.Cells.Item(11, 1).Value = $newRow11[0]
$ws.Cells.Item(11, 2).Value = $newRow11[1]
$ws.Cells.Item(11, 3).Value = $newRow11[2]
$ws.Cells.Item(11, 4).Value = $newRow11[3]
$ws.Cells.Item(11, 5).Value = $newRow11[4]
$ws.Cells.Item(11, 6).Value = $newRow11[5]
$ws.Cells.Item(11, 7).Value = $newRow11[6]
$ws.Cells.Item(11, 8).NumberFormat = "@"
$ws.Cells.Item(11, 8).Value = $newRow11[7]
$ws.Cells.Item(11, 9).Value = $newRow11[8]
$ws.Cells.Item(11, 10).Value = $newRow11[9]
$ws.Cells.Item(11, 11).Value = $newRow11[10]
$ws.Cells.Item(11, 12).Value = $newRow11[11]
$ws.Cells.Item(11, 13).Value = $newRow11[12]
$ws.Cells.Item(11, 14).Value = $newRow11[13]

# Row 12 --------------------------------------------------------------
$ws.Cells.Item(12, 1).NumberFormat = "@"
$ws.Cells.Item(12, 1).Value = $newRow12[0]
$ws.Cells.Item(12, 2).Value = $newRow12[1]
$ws.Cells.Item(12, 3).Value = $newRow12[2]
$ws.Cells.Item(12, 4).Value = $newRow12[3]
$ws.Cells.Item(12, 5).Value = $newRow12[4]
$ws.Cells.Item(12, 6).Value = $newRow12[5]
$ws.Cells.Item(12, 7).Value = $newRow12[6]
$ws.Cells.Item(12, 8).NumberFormat = "@"
$ws.Cells.Item(12, 8).Value = $newRow12[7]
$ws.Cells.Item(12, 9).Value = $newRow12[8]
$ws.Cells.Item(12, 10).Value = $newRow12[9]
$ws.Cells.Item(12, 11).Value = $newRow12[10]
$ws.Cells.Item(12, 12).Value = $newRow12[11]
$ws.Cells.Item(12, 13).Value = $newRow12[12]
$ws.Cells.Item(12, 14).Value = $newRow12[13]

# Row 13 --------------------------------------------------------------
$ws.Cells.Item(13, 1).NumberFormat = "@"
$ws.Cells.Item(13, 1).Value = $newRow13[0]
$ws.Cells.Item(13, 2).Value = $newRow13[1]
$ws.Cells.Item(13, 3).Value = $newRow13[2]
$ws.Cells.Item(13, 4).Value = $newRow13[3]
$ws.Cells.Item(13, 5).Value = $newRow13[4]
$ws.Cells.Item(13, 6).Value = $newRow13[5]
$ws.Cells.Item(13, 7).Value = $newRow13[6]
$ws.Cells.Item(13, 8).NumberFormat = "@"
$ws.Cells.Item(13, 8).Value = $newRow13[7]
$ws.Cells.Item(13, 9).Value = $newRow13[8]
$ws.Cells.Item(13, 10).Value = $newRow13[9]
$ws.Cells.Item(13, 11).NumberFormat = "@"
$ws.Cells.Item(13, 11).Value = $newRow13[10]
$ws.Cells.Item(13, 12).Value = $newRow13[11]
$ws.Cells.Item(13, 13).Value = $newRow13[12]
$ws.Cells.Item(13, 14).Value = $newRow13[13]
